$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.880.81"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "1.638.51"
$ws.Range("E3").Value = "  -0.86%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "215.62"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.5033"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.60%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2567"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06391"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.48%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "19.72"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.40%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07721"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.651.94"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.279"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "1.868.31"
$ws.Range("E14").Value = "  -0.62%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.5464"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "0.0₅7915"
$ws.Range("E16").Value = "  -1.41%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "64.01"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "25.917.99"
$ws.Range("E18").Value = "  -0.89%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "202.03"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.81%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.369"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.85%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "9.930"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.46%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.982"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.16%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.04%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.923"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +10.15%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "141.86"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.34%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.1139"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.51%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.68"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.82%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "6.722"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -3.81%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.244"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.13%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.04991"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.76%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.277"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.04%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.199"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.62%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.542"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.48%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.375"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("D36").Value = "1.168.58"
$ws.Range("E36").Value = "  +0.28%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.631"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.19%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.8920"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.67%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.5592"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.75%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.01564"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.34%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.694"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.68%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.8071"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.68%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "99.80"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  -5.91%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.4525"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("E48").Value = "  +0.08%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "54.95"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.12%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.05075"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.35%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
